# Apply the "Price Option Page check for hints regarding mandatory fields"
# test row to the automation checklist sheet (msz - video 5 including today keyword)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 43: Action / Modus / ... / Record-Control / Expected result
$ws.Range("A43").Value = "Price Option Page check for hints regarding mandatory fields"
$ws.Range("B43").Value = "<CHK>"
$ws.Range("F43").Value = "Price Option Page check for hints regarding mandatory fields"
$ws.Range("H43").Value = "<NOP>"

# Column F widened to fit the new, longer text (closest value reachable
# through the exposed ColumnWidth property, which this runtime snaps to
# 1/6-character increments).
$ws.Columns.Item(6).ColumnWidth = 49.6666666666667

# Move the view / selection down to the newly added row, matching the
# author's on-screen focus after adding the row.
$null = $ws.Range("A43").Select()
